# Update column G (K = strikeouts) values for rows 2-41 on the active sheet.
# The commit message explains that the save-data generation script now
# writes actual strikeout counts (K) instead of the previous "Strike#"
# proxy values, so these are literal data updates row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 3
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 0
    17 = 2
    18 = 2
    19 = 2
    20 = 3
    21 = 3
    22 = 1
    23 = 2
    24 = 1
    25 = 2
    26 = 2
    27 = 1
    28 = 2
    29 = 5
    30 = 0
    31 = 2
    32 = 3
    33 = 2
    34 = 1
    35 = 3
    36 = 1
    37 = 0
    38 = 0
    39 = 0
    40 = 1
    41 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
